$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New issue row (row 4): "Difficolta nel rappresentare l'avanzamento del
# torneo con l'ERM" -- filled in first so its shared strings land before the
# row 2 strings (matches the authored sharedStrings.xml ordering).
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Difficoltà nel rappresentare l'avanzamento del torneo con l'ERM"
$ws.Range("C4").Value = "Massimo Rizzoli"
$ws.Range("D4").Value = 43158
$ws.Range("E4").Value = "Dimitri Malferrari, Massimo Rizzoli, Luca Milano"
$ws.Range("F4").Value = "Utilizzo di id numerici per partite e squadre univoci all'interno di un torneo"
$ws.Range("G4").Value = "Completato"

# --- New issue row (row 2): "Difficolta nel caricare cartelle vuote sulla
# repository github"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Difficoltà nel caricare cartelle vuote sulla repository github"
$ws.Range("C2").Value = "Dimitri Malferrari"
$ws.Range("D2").Value = 43139
$ws.Range("E2").Value = "Dimitri Malferrari, Massimo Rizzoli, Luca Milano"
$ws.Range("F2").Value = "Aggiungere un file chiamato .keep vuoto"
$ws.Range("G2").Value = "Completato"

# --- Existing issue ("Difficolta nel visualizzare file ERM") moves down to
# row 3, reusing the same text (shared strings already present).
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Difficoltà nel visualizzare file ERM"
$ws.Range("C3").Value = "Tommaso Eccher"
$ws.Range("D3").Value = 43171
$ws.Range("E3").Value = "Filippo Pizzini"
$ws.Range("F3").Value = "Copiare in un file di testo rinominato *.XML"
$ws.Range("G3").Value = "Completato"

# --- Extend the blank/date-formatted placeholder rows by one (row 29),
# copying the formatting of the last existing placeholder row (28).
$ws.Range("D28").Copy($ws.Range("D29"))

# --- Column width adjustments (bestFit-style widening of B, E, F to fit the
# new, longer text).
$ws.Columns("B").ColumnWidth = 58.166666666666664
$ws.Columns("E").ColumnWidth = 43.166666666666664
$ws.Columns("F").ColumnWidth = 68

# --- Selection / view state.
$ws.Range("H2").Select()
